$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.203892230987549
$ws.Range("B1").Value = 1.918632507324219
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.930150508880615
$ws.Range("E1").Value = 1.204926609992981
